$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G1").Value = "voterId"
$ws.Range("H4").Select()
